$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.689.87'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.844.66'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '315.45'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4318'
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("D8").Value = '0.3710'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("D9").Value = '0.07332'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '0.8785'
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").Value = '21.07'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("D12").Value = '1.924.82'
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").Value = '5.478'
$ws.Range("E13").Value = '  +2.83%  '
$ws.Range("D14").Value = '6.608'
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '0.06954'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '81.16'
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").Value = '0.000009057'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '15.57'
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("D21").Value = '28.076.14'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '5.113'
$ws.Range("E22").Value = '  +3.16%  '
$ws.Range("D23").Value = '11.00'
$ws.Range("E23").Value = '  +5.86%  '
$ws.Range("D24").Value = '2.143.04'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '154.10'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").Value = '5.327'
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("D29").Value = '115.87'
$ws.Range("E29").Value = '  -4.59%  '
$ws.Range("D30").Value = '1.878'
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("D31").Value = '0.08917'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  +3.38%  '
$ws.Range("D33").Value = '4.623'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").Value = '1.176'
$ws.Range("E34").Value = '  +6.70%  '
$ws.Range("D35").Value = '2.974'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '0.05445'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("D39").Value = '0.01963'
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").Value = '2.840'
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = '0.5180'
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("E42").Value = '  +2.49%  '
$ws.Range("D43").Value = '6.817'
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").Value = '8.659'
$ws.Range("E44").Value = '  +3.38%  '
$ws.Range("E45").Value = '  +3.42%  '
$ws.Range("D46").Value = '0.4798'
$ws.Range("E46").Value = '  +2.72%  '
$ws.Range("D47").Value = '0.06558'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").Value = '106.55'
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").Value = '1.664'
$ws.Range("E50").Value = '  +2.71%  '
$ws.Range("D51").Value = '1.849'
$ws.Range("E51").Value = '  +6.09%  '
